# corrected filepath for income_counts_entry_and_exit function
#
# The rows of Program/percent-income data were scrambled by the previous
# (incorrect) filepath; this fixes the row order so the labels in column A
# line up with the correct percentages in columns B and C.
#
# Rows 1, 5 and 6 are already correct and are left untouched; only rows
# 2, 3, 4, 7, 8 and 9 need their contents corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row number -> (Program label, % income at entry, % income at exit)
$updates = @{
    2 = @("11495", 10.53, 30.77)
    3 = @("1371",   3.01,  5.21)
    4 = @("143",    1.95,  4.95)
    7 = @("MC",     5.43,  5.69)
    8 = @("OC",     4.02,  4.97)
    9 = @("SPC",    6.58, 23.4)
}

foreach ($r in $updates.Keys) {
    $entry = $updates[$r]

    # Column A labels can look numeric (e.g. "11495"); force text formatting
    # so they stay strings rather than being auto-converted to numbers.
    $ws.Cells.Item($r, 1).NumberFormat = "@"
    $ws.Cells.Item($r, 1).Value = [string]$entry[0]

    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
}
